# Added filtering options for the Component Analysis
# Clears the trailing forecast-horizon error cells that are no longer part
# of the evaluated (filtered) range for each row of the component results
# table, removing them from the sheet entirely (matching the source diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clearAddress = "I2:K2,G3:K3,I4:K4,G5:K5,I6:K6,G7:K7,I8:K8,G9:K9,I10:K10,G11:K11,I12:K12,K13,I14:K14,K15,I16:K16,K17,J18:K18,I19:K19,K21,J22:K22,I23:K23,K25,J26:K26,I27:K27,K29,J30:K30,I31:K31,K33,J34:K34,I35:K35,K37,J38:K38,I39:J39"

# The interop runtime only clears the first area of a multi-area Range, so
# iterate through each disjoint area explicitly to make sure every one of
# them is cleared (matching the diff, which removes many small, separate
# groups of cells across different rows).
$clearRange = $ws.Range($clearAddress)
foreach ($area in $clearRange.Areas) {
    $area.ClearContents()
}
